# Regenerate merged AHB files
# 1. Rename the comparison-column headers in row 1:
#      *_old  -> *_FV2404   (columns A-J)
#      *_new  -> *_FV2410   (columns L-U)
#    Column K ("diff") is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($col = 1; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = ($cell.Text -replace "_old$", "_FV2404")
}
for ($col = 12; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = ($cell.Text -replace "_new$", "_FV2410")
}

# 2. Turn the data range into an Excel Table ("Table1") with autofilter.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U58"), $null, 1)
$tbl.Name = "Table1"

# 3. Freeze the header row (split below row 1, top-left of the scrolling
#    pane is A2).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
